# Sales module update: refresh reference number / dates on the Purchase
# Order header sheet and drop the unused trailing blank rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("PurchaseOrderHeader")

# --- Update the header-row text formats (apply Text number format so the
#     date-like values keep rendering literally). This mirrors selecting
#     the header cells and setting their format to Text.
$ws.Range("C1").NumberFormat = "@"
$ws.Range("E1").NumberFormat = "@"

# --- Data row values
$ws.Range("B2").Value = "231225-1"
$ws.Range("C2").Value = "24-12-2025"
$ws.Range("E2").Value = "30-12-2025"

# --- Remove the extra blank, styled rows 6-14 (no longer needed)
$ws.Range("A6:E14").EntireRow.Delete()

# --- Restore the active selection to E3, matching the saved view state
$ws.Range("E3").Select()
